$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "z-mart"
$ws.Cells.Item(2, 2).Value = 9788
$ws.Cells.Item(2, 3).Value = "Polo"
$ws.Cells.Item(2, 4).Value = 45
$ws.Cells.Item(2, 5).Value = "Kilogram"
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = "2nd Dec, 2022"

# Row 3
$ws.Cells.Item(3, 1).Value = "T-shirts"
$ws.Cells.Item(3, 2).Value = 342
$ws.Cells.Item(3, 3).Value = "Armany"
$ws.Cells.Item(3, 4).Value = 120
$ws.Cells.Item(3, 5).Value = "Piece"
$ws.Cells.Item(3, 6).Value = 6
$ws.Cells.Item(3, 7).Value = "5th Dec, 2022"

# Row 4
$ws.Cells.Item(4, 1).Value = "Product 111"
$ws.Cells.Item(4, 2).Value = 85958586
$ws.Cells.Item(4, 3).Value = "Zara"
$ws.Cells.Item(4, 4).Value = 85
$ws.Cells.Item(4, 5).Value = "Kilogram"
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = "5th Dec, 2022"

# Row 5
$ws.Cells.Item(5, 1).Value = "Product 222"
$ws.Cells.Item(5, 2).Value = 46598685
$ws.Cells.Item(5, 3).Value = "Zara"
$ws.Cells.Item(5, 4).Value = 96
$ws.Cells.Item(5, 5).Value = "Piece"
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = "5th Dec, 2022"

# Row 6
$ws.Cells.Item(6, 1).Value = "Product 333"
$ws.Cells.Item(6, 2).Value = 53258695
$ws.Cells.Item(6, 3).Value = "iPhone"
$ws.Cells.Item(6, 4).Value = 63
$ws.Cells.Item(6, 5).Value = "Meter"
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = "5th Dec, 2022"
